$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the target sentence boundary inside the paragraph that talks about
# Datalab attendance / corona exemptions. We work purely with plain-text
# offsets on $d.Content.Text so the script is resilient to how the engine
# happens to have split the existing runs.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$marker = "of corona, excused from attending because of personal circumstances by staff. "
$idx = $full.IndexOf($marker)
if ($idx -lt 0) {
    throw "Could not locate the target sentence in the document."
}
$p1end = $idx + $marker.Length   # boundary right after "...by staff. "

# ---------------------------------------------------------------------------
# Step 1: split the paragraph right after "...by staff. " so the remaining
# text ("Failing to attend ... suspension.") becomes its own paragraph. This
# is a boundary-safe operation that cleanly separates the runs.
# ---------------------------------------------------------------------------
$rSplit1 = $d.Range($p1end, $p1end)
$rSplit1.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# Step 2: insert the new sentence at the start of what is now the second
# paragraph (right after the paragraph mark we just created).
# ---------------------------------------------------------------------------
$p2start = $p1end + 1
$newSentence = "Being late twice for more than 15 minutes without being excused also constitutes non-attendance. "
$rInsert = $d.Range($p2start, $p2start)
$rInsert.InsertBefore($newSentence)

# ---------------------------------------------------------------------------
# Step 3: split again right after the newly inserted sentence so it becomes
# an isolated paragraph of its own - this lets us apply character formatting
# (LanguageID) precisely to only this sentence.
# ---------------------------------------------------------------------------
$p2sentenceEnd = $p2start + $newSentence.Length
$rSplit2 = $d.Range($p2sentenceEnd, $p2sentenceEnd)
$rSplit2.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# Step 4: apply the "en-NL" language mark to the whole (now-isolated)
# paragraph containing only the new sentence. This produces
# <w:lang w:val="en-NL" w:eastAsia="en-NL"/> on that run, matching the
# target formatting for the newly-added sentence.
# ---------------------------------------------------------------------------
$paragraphs = $d.Paragraphs
$count = $paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $paragraphs.Item($i)
    if ($para.Range.Start -le $p2start -and $para.Range.End -gt $p2start) {
        $para.Range.LanguageID = "en-NL"
        break
    }
}

# ---------------------------------------------------------------------------
# Step 5: merge the three paragraphs back into a single paragraph by
# deleting the two paragraph marks we introduced. Delete the later mark
# first so the earlier offset stays valid.
# ---------------------------------------------------------------------------
$mark2 = $d.Range($p2sentenceEnd, $p2sentenceEnd + 1)
$mark2.Delete()

$mark1 = $d.Range($p1end, $p1end + 1)
$mark1.Delete()

# ---------------------------------------------------------------------------
# Step 6: remove the now-stale "next year" qualifier from the retake
# sentence: "...retake it in its entirety next year;" -> "...retake it in
# its entirety;". This edit is fully contained within the trailing run, so
# it does not disturb the runs we just built above.
# ---------------------------------------------------------------------------
$full2 = $d.Content.Text
$oldTail = "retake it in its entirety next year;"
$newTail = "retake it in its entirety;"
$tailIdx = $full2.IndexOf($oldTail)
if ($tailIdx -lt 0) {
    throw "Could not locate the 'next year' text to remove."
}
$rTail = $d.Range($tailIdx, $tailIdx + $oldTail.Length)
$rTail.Text = $newTail
